# Scheduled runner update: refresh market-price-derived columns (H-N) on each
# job sheet (currentAveragePrice.. LeveProfitHQ) with newly pulled pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: The Writing Is Not on the Wall / Enchanted Silver Ink
$ws.Range("H28").Value = 6666853
$ws.Range("I28").Value = 140.0625
$ws.Range("J28").Value = 18518788
$ws.Range("K28").Value = 140.0625
$ws.Range("L28").Value = 18518788
$ws.Range("M28").Value = 344.9375
$ws.Range("N28").Value = -18519758

# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 807.1818
$ws.Range("I62").Value = 775.44446
$ws.Range("K62").Value = 775.44446
$ws.Range("M62").Value = -151.44446

# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 807.1818
$ws.Range("I65").Value = 775.44446
$ws.Range("K65").Value = 3877.2223
$ws.Range("M65").Value = -757.2223000000004

# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 4520.909
$ws.Range("I76").Value = 3361.4285
$ws.Range("K76").Value = 3361.4285
$ws.Range("M76").Value = -3046.4285

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 4520.909
$ws.Range("I79").Value = 3361.4285
$ws.Range("K79").Value = 3361.4285
$ws.Range("M79").Value = -2269.4285

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 2082
$ws.Range("I98").Value = 2140.8
$ws.Range("K98").Value = 2140.8
$ws.Range("M98").Value = -642.8000000000002

# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 48726.094
$ws.Range("I100").Value = 56463.777
$ws.Range("J100").Value = 2300
$ws.Range("K100").Value = 56463.777
$ws.Range("L100").Value = 2300
$ws.Range("M100").Value = -55922.777
$ws.Range("N100").Value = -3382

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 2082
$ws.Range("I122").Value = 2140.8
$ws.Range("K122").Value = 6422.400000000001
$ws.Range("M122").Value = -3972.400000000001

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1591.1569
$ws.Range("I138").Value = 937.2381
$ws.Range("J138").Value = 2048.9
$ws.Range("K138").Value = 2811.7143
$ws.Range("L138").Value = 6146.700000000001
$ws.Range("M138").Value = 2328.2857
$ws.Range("N138").Value = -16426.7

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 855.3889
$ws.Range("I2").Value = 641.5714
$ws.Range("J2").Value = 991.4545000000001
$ws.Range("K2").Value = 641.5714
$ws.Range("L2").Value = 991.4545000000001
$ws.Range("M2").Value = -528.5714
$ws.Range("N2").Value = -1217.4545

# Row 88: The Mast Chance / Adamantite Rivets
$ws.Range("H88").Value = 46671.332
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 46671.332
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 46671.332
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -47483.332

# Row 91: The Rose and the Riveter (L) / Adamantite Rivets
$ws.Range("H91").Value = 46671.332
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 46671.332
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 46671.332
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -49479.332

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 855.3889
$ws.Range("I116").Value = 641.5714
$ws.Range("J116").Value = 991.4545000000001
$ws.Range("K116").Value = 641.5714
$ws.Range("L116").Value = 991.4545000000001
$ws.Range("M116").Value = 1652.4286
$ws.Range("N116").Value = -5579.4545

# Row 117: Signed, Shield, Delivered / Titanbronze Tower Shield
$ws.Range("H117").Value = 27115.334
$ws.Range("J117").Value = 27115.334
$ws.Range("L117").Value = 27115.334
$ws.Range("N117").Value = -36293.334

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 855.3889
$ws.Range("I3").Value = 641.5714
$ws.Range("J3").Value = 991.4545000000001
$ws.Range("K3").Value = 641.5714
$ws.Range("L3").Value = 991.4545000000001
$ws.Range("M3").Value = -527.5714
$ws.Range("N3").Value = -1219.4545

# Row 10: Bring Me the Head Knife of Al'bedo Derssia / Bronze Head Knife
$ws.Range("H10").Value = 315.2
$ws.Range("I10").Value = 244
$ws.Range("J10").Value = 600
$ws.Range("K10").Value = 244
$ws.Range("L10").Value = 600
$ws.Range("M10").Value = -104
$ws.Range("N10").Value = -880

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 1648.125
$ws.Range("J86").Value = 1335.6666
$ws.Range("L86").Value = 1335.6666
$ws.Range("N86").Value = -3581.6666

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 1648.125
$ws.Range("J89").Value = 1335.6666
$ws.Range("L89").Value = 6678.333000000001
$ws.Range("N89").Value = -17910.333

# Row 137: Dagger Swagger / Cobalt Tungsten Khukuri
$ws.Range("H137").Value = 39180
$ws.Range("J137").Value = 39180
$ws.Range("L137").Value = 39180
$ws.Range("N137").Value = -49380

# Row 138: Bladewinner / Titanium Gold Greatsword
$ws.Range("H138").Value = 22877
$ws.Range("J138").Value = 22877
$ws.Range("L138").Value = 22877
$ws.Range("N138").Value = -33157

# Row 140: Ceremonial Teeth / Ra'Kaznar Twinfangs
$ws.Range("H140").Value = 33000
$ws.Range("J140").Value = 33000
$ws.Range("L140").Value = 33000
$ws.Range("N140").Value = -43360

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 14471.714
$ws.Range("I31").Value = 1684.65
$ws.Range("J31").Value = 46439.375
$ws.Range("K31").Value = 1684.65
$ws.Range("L31").Value = 46439.375
$ws.Range("M31").Value = -1389.65
$ws.Range("N31").Value = -47029.375

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 14471.714
$ws.Range("I34").Value = 1684.65
$ws.Range("J34").Value = 46439.375
$ws.Range("K34").Value = 1684.65
$ws.Range("L34").Value = 46439.375
$ws.Range("M34").Value = -1482.65
$ws.Range("N34").Value = -46843.375

$ws = $wb.Worksheets.Item("CUL")
# Row 59: Comfort Me with Mushrooms / Buttons in a Blanket
$ws.Range("H59").Value = 1498.875
$ws.Range("I59").Value = 700.3333
$ws.Range("J59").Value = 1978
$ws.Range("K59").Value = 2100.9999
$ws.Range("L59").Value = 5934
$ws.Range("M59").Value = -1560.9999
$ws.Range("N59").Value = -7014

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 3396.25
$ws.Range("I80").Value = 3005
$ws.Range("J80").Value = 3526.6667
$ws.Range("K80").Value = 3005
$ws.Range("L80").Value = 3526.6667
$ws.Range("M80").Value = -2007
$ws.Range("N80").Value = -5522.6667

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 3396.25
$ws.Range("I83").Value = 3005
$ws.Range("J83").Value = 3526.6667
$ws.Range("K83").Value = 15025
$ws.Range("L83").Value = 17633.3335
$ws.Range("M83").Value = -10033
$ws.Range("N83").Value = -27617.3335

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 1286.6522
$ws.Range("I122").Value = 946.2143
$ws.Range("J122").Value = 1816.2222
$ws.Range("K122").Value = 2838.6429
$ws.Range("L122").Value = 5448.6666
$ws.Range("M122").Value = -388.6428999999998
$ws.Range("N122").Value = -10348.6666

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 1482118.1
$ws.Range("I7").Value = 2115
$ws.Range("J7").Value = 4277679.5
$ws.Range("K7").Value = 2115
$ws.Range("L7").Value = 4277679.5
$ws.Range("M7").Value = -2003
$ws.Range("N7").Value = -4277903.5

# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 2915.4707
$ws.Range("I40").Value = 2437.5334
$ws.Range("J40").Value = 6500
$ws.Range("K40").Value = 2437.5334
$ws.Range("L40").Value = 6500
$ws.Range("M40").Value = -2301.5334
$ws.Range("N40").Value = -6772

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 1482118.1
$ws.Range("I126").Value = 2115
$ws.Range("J126").Value = 4277679.5
$ws.Range("K126").Value = 6345
$ws.Range("L126").Value = 12833038.5
$ws.Range("M126").Value = -3875
$ws.Range("N126").Value = -12837978.5
